$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new text value for D (Price) and E (Volume 1h) columns
# D-column cells are forced to Text format to preserve exact textual representation
# (avoids Excel auto-converting strings like '1.000' or '0.09990' into numbers and
# losing significant trailing zeros / multi-dot grouping).
$updates = [ordered]@{
    D2 = '30.631.68'
    D3 = '1.872.14'
    E3 = '  +0.00%  '
    D4 = '1.000'
    E4 = '  +0.01%  '
    D5 = '248.14'
    E5 = '  +1.22%  '
    D6 = '0.9997'
    E6 = '  -0.04%  '
    D7 = '0.4734'
    E7 = '  +0.30%  '
    D8 = '0.2913'
    E8 = '  +1.28%  '
    D9 = '0.06476'
    E9 = '  +0.06%  '
    D10 = '22.05'
    E10 = '  +4.61%  '
    D11 = '0.07693'
    E11 = '  -1.08%  '
    D12 = '96.57'
    E12 = '  +1.38%  '
    D13 = '0.7376'
    E13 = '  +2.94%  '
    D14 = '1.869.93'
    E14 = '  -0.11%  '
    E15 = '  +0.71%  '
    D16 = '272.67'
    E16 = '  -1.56%  '
    D17 = '30.646.74'
    E17 = '  +1.03%  '
    D18 = '13.33'
    E18 = '  -0.34%  '
    D19 = '0.9992'
    E19 = '  -0.07%  '
    D20 = '0.000007513'
    E20 = '  -0.83%  '
    D21 = '2.112.57'
    E21 = '  -0.13%  '
    D22 = '0.9984'
    E22 = '  -0.15%  '
    D23 = '5.263'
    E23 = '  +0.44%  '
    D24 = '6.176'
    D25 = '9.218'
    E25 = '  -0.53%  '
    D26 = '163.92'
    E26 = '  -1.17%  '
    D27 = '18.77'
    E27 = '  -0.75%  '
    D28 = '1.910'
    E28 = '  -0.16%  '
    D29 = '0.09990'
    E29 = '  +0.93%  '
    E30 = '  -2.59%  '
    D31 = '1.511'
    E31 = '  -0.34%  '
    D32 = '4.280'
    E32 = '  -0.03%  '
    D33 = '4.101'
    E33 = '  +1.70%  '
    D34 = '0.04795'
    E34 = '  +0.33%  '
    E35 = '  -0.36%  '
    D36 = '0.6963'
    E36 = '  -0.02%  '
    D37 = '2.714'
    E37 = '  -0.23%  '
    D38 = '0.01850'
    E38 = '  -0.17%  '
    D39 = '2.755'
    E39 = '  +0.47%  '
    D40 = '6.212'
    E40 = '  -2.67%  '
    D41 = '73.19'
    E41 = '  +4.06%  '
    D42 = '1.970'
    E42 = '  +2.64%  '
    D43 = '0.4181'
    E43 = '  +1.52%  '
    D44 = '0.9996'
    E44 = '  -0.03%  '
    E45 = '  -1.25%  '
    D46 = '101.82'
    E46 = '  -0.26%  '
    D47 = '9.312'
    E47 = '  +0.19%  '
    E48 = '  +0.37%  '
    D49 = '6.971'
    E49 = '  -1.80%  '
    D50 = '917.63'
    E50 = '  -0.34%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref.StartsWith('D')) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}
